$wb = $excel.ActiveWorkbook

# This script applies refreshed market-data values (columns H-N) produced by
# the scheduled Jenova_Profits market data runner across all 8 job sheets.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 963.4651
$ws.Range("J17").Value = 963.4651
$ws.Range("L17").Value = 2890.3953
$ws.Range("N17").Value = -3226.3953
$ws.Range("H34").Value = 15500
$ws.Range("I34").Value = 15500
$ws.Range("K34").Value = 15500
$ws.Range("M34").Value = -15297
$ws.Range("H36").Value = 15500
$ws.Range("I36").Value = 15500
$ws.Range("K36").Value = 15500
$ws.Range("M36").Value = -14785
$ws.Range("H43").Value = 1087.6
$ws.Range("J43").Value = 1087.6
$ws.Range("L43").Value = 1087.6
$ws.Range("N43").Value = -1225.6
$ws.Range("H55").Value = 101405.3
$ws.Range("I55").Value = 490
$ws.Range("J55").Value = 112618.11
$ws.Range("K55").Value = 490
$ws.Range("L55").Value = 112618.11
$ws.Range("M55").Value = -276
$ws.Range("N55").Value = -113046.11
$ws.Range("H107").Value = 85354.914
$ws.Range("I107").Value = 101906.3
$ws.Range("J107").Value = 2598
$ws.Range("K107").Value = 101906.3
$ws.Range("L107").Value = 2598
$ws.Range("M107").Value = -99986.3
$ws.Range("N107").Value = -6438
$ws.Range("H137").Value = 2547.535
$ws.Range("I137").Value = 1363.8
$ws.Range("J137").Value = 5279.231
$ws.Range("K137").Value = 4091.4
$ws.Range("L137").Value = 15837.693
$ws.Range("M137").Value = -1541.4
$ws.Range("N137").Value = -20937.693
$ws.Range("H138").Value = 5612.9
$ws.Range("I138").Value = 4107.636
$ws.Range("J138").Value = 6037.4614
$ws.Range("K138").Value = 12322.908
$ws.Range("L138").Value = 18112.3842
$ws.Range("M138").Value = -7182.908000000001
$ws.Range("N138").Value = -28392.3842

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6766
$ws.Range("I45").Value = 2899
$ws.Range("K45").Value = 2899
$ws.Range("M45").Value = -2522
$ws.Range("H122").Value = 6202.2
$ws.Range("I122").Value = 2012
$ws.Range("K122").Value = 6036
$ws.Range("M122").Value = -3586
$ws.Range("H132").Value = 5255.961
$ws.Range("I132").Value = 2848.1353
$ws.Range("K132").Value = 8544.4059
$ws.Range("M132").Value = -6014.4059
$ws.Range("H141").Value = 74333.336
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 96500
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 96500
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -106860

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -1226
$ws.Range("H3").Value = 23597.6
$ws.Range("J3").Value = 4664.3335
$ws.Range("L3").Value = 4664.3335
$ws.Range("N3").Value = -4890.3335
$ws.Range("H31").Value = 3302.8
$ws.Range("I31").Value = 1224.6
$ws.Range("J31").Value = 7459.2
$ws.Range("K31").Value = 1224.6
$ws.Range("L31").Value = 7459.2
$ws.Range("M31").Value = -929.5999999999999
$ws.Range("N31").Value = -8049.2
$ws.Range("H34").Value = 3302.8
$ws.Range("I34").Value = 1224.6
$ws.Range("J34").Value = 7459.2
$ws.Range("K34").Value = 1224.6
$ws.Range("L34").Value = 7459.2
$ws.Range("M34").Value = -1022.6
$ws.Range("N34").Value = -7863.2
$ws.Range("H43").Value = 21578.5
$ws.Range("J43").Value = 21578.5
$ws.Range("L43").Value = 21578.5
$ws.Range("N43").Value = -21946.5
$ws.Range("H62").Value = 3511
$ws.Range("I62").Value = 3499.8333
$ws.Range("K62").Value = 3499.8333
$ws.Range("M62").Value = -2875.8333
$ws.Range("H65").Value = 3511
$ws.Range("I65").Value = 3499.8333
$ws.Range("K65").Value = 17499.1665
$ws.Range("M65").Value = -14379.1665
$ws.Range("H69").Value = 29999.666
$ws.Range("I69").Value = 29999.666
$ws.Range("K69").Value = 29999.666
$ws.Range("M69").Value = -29250.666
$ws.Range("H72").Value = 29999.666
$ws.Range("I72").Value = 29999.666
$ws.Range("K72").Value = 89998.99800000001
$ws.Range("M72").Value = -86254.99800000001
$ws.Range("H94").Value = 2924.6
$ws.Range("J94").Value = 3570.3333
$ws.Range("L94").Value = 3570.3333
$ws.Range("N94").Value = -4472.3333
$ws.Range("H101").Value = 21578.5
$ws.Range("J101").Value = 21578.5
$ws.Range("L101").Value = 21578.5
$ws.Range("N101").Value = -28068.5
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72080
$ws.Range("H132").Value = 4271.15
$ws.Range("J132").Value = 6030
$ws.Range("L132").Value = 18090
$ws.Range("N132").Value = -23150

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.26087
$ws.Range("J2").Value = 246.11111
$ws.Range("L2").Value = 1476.66666
$ws.Range("N2").Value = -1702.66666
$ws.Range("H12").Value = 297.6
$ws.Range("J12").Value = 332.07693
$ws.Range("L12").Value = 996.2307900000001
$ws.Range("N12").Value = -1342.23079
$ws.Range("H14").Value = 1070.25
$ws.Range("I14").Value = 1070.25
$ws.Range("K14").Value = 3210.75
$ws.Range("M14").Value = -3037.75
$ws.Range("H92").Value = 891
$ws.Range("I92").Value = 638.8
$ws.Range("J92").Value = 1031.1111
$ws.Range("K92").Value = 1916.4
$ws.Range("L92").Value = 3093.3333
$ws.Range("M92").Value = -668.3999999999999
$ws.Range("N92").Value = -5589.3333
$ws.Range("H102").Value = 15000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H133").Value = 12244.5
$ws.Range("I133").Value = 21989
$ws.Range("K133").Value = 65967
$ws.Range("M133").Value = -60907
$ws.Range("M102").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 166668540
$ws.Range("I126").Value = 250000940
$ws.Range("J126").Value = 3749.5
$ws.Range("K126").Value = 750002820
$ws.Range("L126").Value = 11248.5
$ws.Range("M126").Value = -750000350
$ws.Range("N126").Value = -16188.5
$ws.Range("H132").Value = 235344.23
$ws.Range("I132").Value = 252545.05
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 757635.1499999999
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -755105.1499999999
$ws.Range("N132").Value = -23060

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3355.8286
$ws.Range("I46").Value = 2935.3333
$ws.Range("K46").Value = 2935.3333
$ws.Range("M46").Value = -2747.3333
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51498
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -157488
$ws.Range("H68").Value = 114684
$ws.Range("I68").Value = 4019.5
$ws.Range("K68").Value = 4019.5
$ws.Range("M68").Value = -3270.5
$ws.Range("H71").Value = 114684
$ws.Range("I71").Value = 4019.5
$ws.Range("K71").Value = 20097.5
$ws.Range("M71").Value = -16353.5
$ws.Range("H132").Value = 5449.2856
$ws.Range("I132").Value = 4249.1665
$ws.Range("J132").Value = 6349.375
$ws.Range("K132").Value = 12747.4995
$ws.Range("L132").Value = 19048.125
$ws.Range("M132").Value = -10217.4995
$ws.Range("N132").Value = -24108.125
$ws.Range("H136").Value = 3376.457
$ws.Range("I136").Value = 3369.9524
$ws.Range("K136").Value = 10109.8572
$ws.Range("M136").Value = -7559.8572

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 128381.5
$ws.Range("J96").Value = 14955
$ws.Range("L96").Value = 14955
$ws.Range("N96").Value = -17701
$ws.Range("H97").Value = 10560
$ws.Range("J97").Value = 10560
$ws.Range("L97").Value = 10560
$ws.Range("N97").Value = -12542
$ws.Range("H98").Value = 227818
$ws.Range("J98").Value = 227818
$ws.Range("L98").Value = 227818
$ws.Range("N98").Value = -233808
$ws.Range("H107").Value = 45446.176
$ws.Range("I107").Value = 54724.367
$ws.Range("J107").Value = 1374.75
$ws.Range("K107").Value = 164173.101
$ws.Range("L107").Value = 4124.25
$ws.Range("M107").Value = -162253.101
$ws.Range("N107").Value = -7964.25
$ws.Range("H126").Value = 2759.5862
$ws.Range("I126").Value = 2769.6843
$ws.Range("J126").Value = 2740.4
$ws.Range("K126").Value = 8309.052899999999
$ws.Range("L126").Value = 8221.200000000001
$ws.Range("M126").Value = -5839.052899999999
$ws.Range("N126").Value = -13161.2
$ws.Range("H132").Value = 5135.409
$ws.Range("I132").Value = 4311.1875
$ws.Range("K132").Value = 12933.5625
$ws.Range("M132").Value = -10403.5625
$ws.Range("H136").Value = 669863
$ws.Range("I136").Value = 717138.9399999999
$ws.Range("K136").Value = 2151416.82
$ws.Range("M136").Value = -2148866.82
